$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-02 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-03 Thursday", 2)

$d.Content.Find.Execute("627×3=", $true, $false, $false, $false, $false, $true, 1, $false, "773×7=", 2)
$d.Content.Find.Execute("125×3=", $true, $false, $false, $false, $false, $true, 1, $false, "776×2=", 2)
$d.Content.Find.Execute("526×8=", $true, $false, $false, $false, $false, $true, 1, $false, "916×4=", 2)
$d.Content.Find.Execute("291×2=", $true, $false, $false, $false, $false, $true, 1, $false, "738×9=", 2)
$d.Content.Find.Execute("335×4=", $true, $false, $false, $false, $false, $true, 1, $false, "762×4=", 2)
$d.Content.Find.Execute("304×7=", $true, $false, $false, $false, $false, $true, 1, $false, "964×2=", 2)
$d.Content.Find.Execute("363×3=", $true, $false, $false, $false, $false, $true, 1, $false, "432×8=", 2)
$d.Content.Find.Execute("607×2=", $true, $false, $false, $false, $false, $true, 1, $false, "184×9=", 2)
$d.Content.Find.Execute("971×5=", $true, $false, $false, $false, $false, $true, 1, $false, "504×8=", 2)
$d.Content.Find.Execute("687×3=", $true, $false, $false, $false, $false, $true, 1, $false, "624×5=", 2)
$d.Content.Find.Execute("364×3=", $true, $false, $false, $false, $false, $true, 1, $false, "695×4=", 2)
$d.Content.Find.Execute("448×8=", $true, $false, $false, $false, $false, $true, 1, $false, "370×7=", 2)
$d.Content.Find.Execute("379×8=", $true, $false, $false, $false, $false, $true, 1, $false, "314×9=", 2)
$d.Content.Find.Execute("499×7=", $true, $false, $false, $false, $false, $true, 1, $false, "780×8=", 2)
$d.Content.Find.Execute("744×7=", $true, $false, $false, $false, $false, $true, 1, $false, "406×4=", 2)
$d.Content.Find.Execute("963×7=", $true, $false, $false, $false, $false, $true, 1, $false, "367×7=", 2)
$d.Content.Find.Execute("852×2=", $true, $false, $false, $false, $false, $true, 1, $false, "278×8=", 2)
$d.Content.Find.Execute("154×5=", $true, $false, $false, $false, $false, $true, 1, $false, "147×6=", 2)
$d.Content.Find.Execute("246×3=", $true, $false, $false, $false, $false, $true, 1, $false, "619×9=", 2)
$d.Content.Find.Execute("949×2=", $true, $false, $false, $false, $false, $true, 1, $false, "563×3=", 2)
$d.Content.Find.Execute("408×7=", $true, $false, $false, $false, $false, $true, 1, $false, "362×9=", 2)
$d.Content.Find.Execute("390×3=", $true, $false, $false, $false, $false, $true, 1, $false, "814×3=", 2)
$d.Content.Find.Execute("669×2=", $true, $false, $false, $false, $false, $true, 1, $false, "603×2=", 2)
$d.Content.Find.Execute("732×2=", $true, $false, $false, $false, $false, $true, 1, $false, "396×2=", 2)
$d.Content.Find.Execute("854×9=", $true, $false, $false, $false, $false, $true, 1, $false, "631×3=", 2)
